$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")

# Add a new user row (row 13) mirroring the structure of the existing
# user rows, representing a newly signed-up user used for the
# "delete user" test setup.
$ws.Range("C13").Value = "helloworld@gmail.com"
$ws.Range("D13").Value = "helloworld"
$ws.Range("E13").Value = "Student"
$ws.Range("F13").Value = $false
